$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 2,3,5,6,8,9 (columns A,B,D,E,F,G,H,P,Q,R,S) were shuffled
# around between each other (a cyclic rotation of full record contents).
# Capture the "before" values of each row first, then write them to their
# new destination rows, so the operation is correct regardless of order.

$cols = @("A","B","D","E","F","G","H","P","Q","R","S")
$rowsInvolved = @(2,3,5,6,8,9)

$snapshot = @{}
foreach ($r in $rowsInvolved) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowData
}

# Mapping: new row (key) receives the old content of row (value)
$mapping = @{
    2 = 8
    3 = 5
    5 = 9
    6 = 2
    8 = 3
    9 = 6
}

foreach ($newRow in $mapping.Keys) {
    $oldRow = $mapping[$newRow]
    $src = $snapshot[$oldRow]
    foreach ($c in $cols) {
        $ws.Range("$c$newRow").Value = $src[$c]
    }
}
